# "Generate Report for Handback"
#
# The localization-status report gains a handback record for
# 790af623-9f78-48ae-afa1-8243b2670b39.md (now "Handed back: in sync with
# en-US", with a Latest Target File / Latest Handback File / Latest Handback
# DateTime filled in on the zh-cn and de-de detail sheets) while
# 05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md drops back to the #2 / row-3 slot
# it vacated ("Ready for handoff", still pending handback).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared URL fragments reused across the hyperlinks we (re)create below.
# ---------------------------------------------------------------------
$url_05a8_md        = "https://github.com/OpenLocalizationTest/oltest/blob/756731c8dfb711571ef2949eef68b03f1c112070/e2e/05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md"
$url_790a_md        = "https://github.com/OpenLocalizationTest/oltest/blob/4bdbe7ee4561750ab8279ab293840c0fffa3982b/e2e/790af623-9f78-48ae-afa1-8243b2670b39.md"
$url_05a8_zhcn_xlf   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f96f76958731798a6dc18875c94281ea2e29c2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.zh-cn.xlf"
$url_790a_zhcn_xlf   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60e84480d896668b102c4fa579619292929b42cb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.zh-cn.xlf"
$url_05a8_dede_xlf   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bbeeb61e811aff4d6c470a6bd4a3b0bf16f28012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.de-de.xlf"
$url_790a_dede_xlf   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1bb660ad387d961644f59102695e528dbbcb06c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.de-de.xlf"

$missing = [System.Type]::Missing

# =======================================================================
# Sheet "Overview" - the file now handed back (790af623...) moves into
# row 2, the file still awaiting handback (05a8c78e...) moves into row 3.
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-29-19 00:29:14"

$wsOverview.Range("A3").Value = "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-28-19 00:28:54"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $url_790a_md, $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $url_05a8_md, $missing, $missing, "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md") | Out-Null

# =======================================================================
# Sheet "zh-cn" - detail rows for the two files, zh-cn target language.
# =======================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2: 790af623... has now been handed back (Include / full round trip).
$wsZhCn.Range("A2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-19 00:29:11"
$wsZhCn.Range("F2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.md"
$wsZhCn.Range("G2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-03-19 00:29:28"
$wsZhCn.Range("I2").Value = "Include"

# Row 3: 05a8c78e... is still only ready for handoff (no handback yet).
$wsZhCn.Range("A3").Value = "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 00:28:52"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url_790a_md,      $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $url_790a_md,      $missing, $missing, ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $url_790a_zhcn_xlf, $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $url_790a_md,      $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $url_790a_zhcn_xlf, $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url_05a8_md,      $missing, $missing, "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $url_05a8_md,      $missing, $missing, ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $url_05a8_zhcn_xlf, $missing, $missing, "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.zh-cn.xlf") | Out-Null

# =======================================================================
# Sheet "de-de" - same shape as "zh-cn" but for the de-de target language.
# =======================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2: 790af623... has now been handed back (Include / full round trip).
$wsDeDe.Range("A2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-19 00:29:14"
$wsDeDe.Range("F2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.md"
$wsDeDe.Range("G2").Value = "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-19 00:29:33"
$wsDeDe.Range("I2").Value = "Include"

# Row 3: 05a8c78e... is still only ready for handoff (no handback yet).
$wsDeDe.Range("A3").Value = "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 00:28:54"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url_790a_md,      $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $url_790a_md,      $missing, $missing, ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $url_790a_dede_xlf, $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $url_790a_md,      $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $url_790a_dede_xlf, $missing, $missing, "790af623-9f78-48ae-afa1-8243b2670b39.929dc35992db41f880409269abbfcf12085f3d6d.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url_05a8_md,      $missing, $missing, "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $url_05a8_md,      $missing, $missing, ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $url_05a8_dede_xlf, $missing, $missing, "05a8c78e-a2c3-4cf0-87b0-2a1dfce82497.0fa7b23741bc5313a69f26378ea049713c511bfa.de-de.xlf") | Out-Null
